$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) - update F column "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 0
$wsExpo.Range("F3").Value = 118
$wsExpo.Range("F4").Value = 1635
$wsExpo.Range("F6").Value = 0
$wsExpo.Range("F8").Value = 0

# Sheet "全部类型" (all types list) - update F column "想去人数" (want-to-go count) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 118
$wsAll.Range("F4").Value = 0
$wsAll.Range("F5").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F7").Value = 0
$wsAll.Range("F8").Value = 0
